# Slide 11 ("Recommendations"), TextBox 5 (shape 2): the five "Increase life
# satisfaction in 5 ways" sub-bullets are being de-italicised / de-mathified
# (the fancy Unicode math-alphanumeric text is replaced by plain text) and
# their bullet glyph switches from the Wingdings "Ø" to a plain Arial "•",
# matching the other top-level bullets in the same box.
#
# We walk the affected paragraphs from bottom to top so that the character
# offsets computed up-front (on the untouched text) remain valid for the
# paragraphs still to be processed (editing a paragraph only changes the
# length of that paragraph itself, never the paragraphs before it).

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(11)
$shp = $s.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

function Set-BulletToArialDot($rng) {
    $bf = $rng.ParagraphFormat.Bullet
    # Set the font before the character so the serializer emits
    # <a:buFont .../><a:buChar .../> in schema order.
    $bf.Font.Name = "Arial"
    $bf.Character = 8226
}

# --- "	Job Involvement" (was "	𝐽𝑜𝑏 Involvement") ------------------------
$r = $tr.Characters(189, 16)
Set-BulletToArialDot $r
$r.Text = "`tJob Involvement"

# --- "	Work life Balance" (was "	Work life 𝐵𝑎𝑙𝑎𝑛𝑐𝑒") ---------------------
$r = $tr.Characters(170, 18)
Set-BulletToArialDot $r
$r.Text = "`tWork life Balance"

# --- "	Relationship Satisfaction" (was "	Relationship 𝑆𝑎𝑡𝑖𝑠faction ") ----
$r = $tr.Characters(142, 27)
Set-BulletToArialDot $r
$r.Text = "`tRelationship Satisfaction"

# --- "	Job Satisfaction" (was "	𝐽𝑜𝑏 𝑆𝑎𝑡𝑖𝑠𝑓𝑎𝑐𝑡𝑖𝑜𝑛 ") ------------------------
$r = $tr.Characters(123, 18)
Set-BulletToArialDot $r
$r.Text = "`tJob Satisfaction"

# --- "	Environmental Satisfaction" -----------------------------------------
# This paragraph is special: the leading tab lives in its own (italic) run,
# separate from the run holding the math-alphanumeric text. Delete the
# separate tab run first so the paragraph is left with a single run (whose
# rPr has no italic), then retype that run's text with the tab restored.
$r = $tr.Characters(97, 25)
Set-BulletToArialDot $r
$tabRun = $tr.Characters(97, 1)
$tabRun.Delete()
$rest = $tr.Characters(97, 24)
$rest.Text = "`tEnvironmental Satisfaction"
